$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" footer text (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 20:35"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1418591
$ws.Range("C4").Value = 9955
$ws.Range("E4").Value = 1033954
$ws.Range("G4").Value = 773
$ws.Range("H4").Value = 84198

# Francia (row 10)
$ws.Range("D10").Value = 58673
$ws.Range("E10").Value = 92478
$ws.Range("F10").Value = 2428
$ws.Range("G10").Value = 83
$ws.Range("H10").Value = 27074

# India (row 15)
$ws.Range("B15").Value = 78042
$ws.Range("C15").Value = 3750
$ws.Range("E15").Value = 49099

# Emiratos Arabes Unidos (row 32)
$ws.Range("B32").Value = 20386
$ws.Range("C32").Value = 725
$ws.Range("D32").Value = 6523
$ws.Range("E32").Value = 13657
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 206

# Belice / Nueva Caledonia swap places (rows 193 & 194 fully swap content)
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0

$ws.Range("A194").Value = "Belice"
$ws.Range("D194").Value = 16
$ws.Range("H194").Value = 2
